$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "Burri"
$wsCuenta.Range("B2").Value = "Pablo Martin"
$wsCuenta.Range("C2").Value = 28263674
$wsCuenta.Range("H14").Select()

# --- DatosHogar sheet (physical file holds the "DatosMotor" NvoNro data due to
#     this workbook's pre-existing sheet/rId relationship scramble) ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 625

# --- DatosMotor sheet (physical file holds the "DatosHogar" SMA00x data, see above) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMA006"
$wsMotor.Range("B2").Value = "ABC12SSMA006"
$wsMotor.Range("C2").Value = "ZAZ123SSMA006"

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200105
$wsAP.Range("A3").Select()
